# Add two new year columns (2021, 2022) to the "17.1.1" revenue table
# and update several historical data points in the existing columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new year headers ---------------------------------------------
$ws.Range("R4").Value = 2021
$ws.Range("S4").Value = 2022

# --- Row 5: "Revenues, total" --------------------------------------------
$ws.Range("P5").Value = 25.6
$ws.Range("Q5").Value = 23.8
$ws.Range("R5").Value = 26.8
$ws.Range("S5").Value = 26.8

# --- Row 6: "Tax revenues" ------------------------------------------------
$ws.Range("P6").Value = 18.6
$ws.Range("Q6").Value = 16.7
$ws.Range("R6").Value = 19.3
$ws.Range("S6").Value = 19.3

# --- Row 7: "Contributions / deductions for social needs" -----------------
$ws.Range("R7").Value = "-"
$ws.Range("S7").Value = "-"

# --- Row 8: "Received official transfers" ---------------------------------
$ws.Range("P8").Value = 2.1
$ws.Range("Q8").Value = 1.8
$ws.Range("R8").Value = 1.8
$ws.Range("S8").Value = 1.8

# --- Row 9: "Non-tax revenues" --------------------------------------------
$ws.Range("P9").Value = 4.9
$ws.Range("Q9").Value = 5.2
$ws.Range("R9").Value = 5.7
$ws.Range("S9").Value = 5.7

# --- Row 10: "Revenues from the sale of non-financial assets" -------------
$ws.Range("R10").Value = 0
$ws.Range("S10").Value = 0

# --- Copy styles from column Q into the new columns R and S ---------------
$ws.Range("Q4:Q10").Copy() | Out-Null
$ws.Range("R4:R10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("Q4:Q10").Copy() | Out-Null
$ws.Range("S4:S10").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$excel.CutCopyMode = 0

# --- Update sheet view selection, matching the authored workbook ----------
$ws.Range("T3").Select() | Out-Null
